$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set B-column grade values for rows that were previously blank
$rowValues = @{
    26 = 2
    27 = 2
    28 = 2
    29 = 2
    30 = 2
    44 = 40
    59 = 5
    60 = 5
    61 = 5
    62 = 5
    63 = 5
    71 = 5
    72 = 5
    73 = 2
    81 = 10
    82 = 10
    83 = 10
    94 = 5
    95 = 5
    105 = 10
    106 = 10
    107 = 10
    108 = 10
    120 = 5
    121 = 5
    122 = 5
    123 = 5
    134 = 0
    135 = 0
    136 = 0
    137 = 0
    138 = 0
    139 = 0
    140 = 0
    141 = 0
    142 = 0
    143 = 0
    144 = 0
    145 = 0
    146 = 0
    147 = 0
    148 = 0
    149 = 0
    150 = 0
    151 = 0
    157 = 0
    167 = 0
    168 = 0
    169 = 0
    170 = 0
    171 = 0
    172 = 0
    173 = 0
}

# Rows whose row height needs to change from 15.75 to 15
$htRows = @(26, 27, 28, 29, 30, 44, 59, 60, 61, 62, 63, 71, 72, 73, 81, 82, 83, 94, 95, 105, 106, 107, 108)

foreach ($r in $rowValues.Keys) {
    $ws.Range("B$r").Value = $rowValues[$r]
}

foreach ($r in $htRows) {
    $ws.Rows($r).RowHeight = 15
}

# Update sheet view: scroll the view to the top, then leave the active
# cell/selection on B171 (matches the saved view state in the workbook)
$ws.Range("A1").Select() | Out-Null
$ws.Range("B171").Select() | Out-Null

Write-Host "Edit complete"
